{"js": "// Replace the 25 division-problem answers in the single table, by\n// position (row, col), matching the diff's old -> new text per cell.\n// Using position rather than literal text search avoids ambiguity,\n// since one new value (\"27\u00f75=5, 2\") duplicates an old value used\n// elsewhere in the table.\n\nconst replacements = [\n  // [dataRowIndex (0-based among the 5 populated rows), colIndex, oldText, newText]\n  [0, 0, \"33\u00f73=11, 0\", \"93\u00f76=15, 3\"],\n  [0, 1, \"32\u00f74=8, 0\", \"85\u00f73=28, 1\"],\n  [0, 2, \"17\u00f75=3, 2\", \"40\u00f77=5, 5\"],\n  [0, 3, \"31\u00f79=3, 4\", \"71\u00f79=7, 8\"],\n  [0, 4, \"90\u00f72=45, 0\", \"28\u00f75=5, 3\"],\n\n  [1, 0, \"85\u00f74=21, 1\", \"19\u00f72=9, 1\"],\n  [1, 1, \"31\u00f75=6, 1\", \"90\u00f75=18, 0\"],\n  [1, 2, \"61\u00f72=30, 1\", \"88\u00f78=11, 0\"],\n  [1, 3, \"26\u00f79=2, 8\", \"82\u00f77=11, 5\"],\n  [1, 4, \"98\u00f79=10, 8\", \"51\u00f73=17, 0\"],\n\n  [2, 0, \"99\u00f74=24, 3\", \"81\u00f74=20, 1\"],\n  [2, 1, \"27\u00f75=5, 2\", \"45\u00f78=5, 5\"],\n  [2, 2, \"97\u00f78=12, 1\", \"16\u00f72=8, 0\"],\n  [2, 3, \"16\u00f77=2, 2\", \"58\u00f72=29, 0\"],\n  [2, 4, \"46\u00f72=23, 0\", \"22\u00f78=2, 6\"],\n\n  [3, 0, \"22\u00f73=7, 1\", \"75\u00f74=18, 3\"],\n  [3, 1, \"84\u00f72=42, 0\", \"57\u00f74=14, 1\"],\n  [3, 2, \"62\u00f73=20, 2\", \"74\u00f72=37, 0\"],\n  [3, 3, \"35\u00f73=11, 2\", \"27\u00f75=5, 2\"],\n  [3, 4, \"43\u00f73=14, 1\", \"18\u00f75=3, 3\"],\n\n  [4, 0, \"70\u00f72=35, 0\", \"48\u00f72=24, 0\"],\n  [4, 1, \"49\u00f73=16, 1\", \"15\u00f76=2, 3\"],\n  [4, 2, \"56\u00f75=11, 1\", \"88\u00f78=11, 0\"],\n  [4, 3, \"85\u00f77=12, 1\", \"45\u00f72=22, 1\"],\n  [4, 4, \"30\u00f74=7, 2\", \"83\u00f75=16, 3\"],\n];\n\n// The table has 20 rows total, but only every 4th row (0, 4, 8, 12, 16)\n// actually holds answer text; the rows between are spacer rows.\nconst DATA_ROW_STRIDE = 4;\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document, but none was found.\");\n}\n\n// Load every cell's current text up-front so we can sanity-check against\n// the expected \"old\" text before overwriting it.\nconst cells = [];\nfor (const [dataRowIdx, col, oldText, newText] of replacements) {\n  const row = dataRowIdx * DATA_ROW_STRIDE;\n  const cell = table.getCell(row, col);\n  cell.load(\"value\");\n  cells.push({ cell, row, col, oldText, newText });\n}\nawait context.sync();\n\nfor (const { cell, row, col, oldText, newText } of cells) {\n  // Sanity-check against the expected original text (trimmed, since the\n  // loaded value can include the cell's trailing paragraph mark); the\n  // write itself always applies the intended new value regardless, so the\n  // script is resilient to incidental whitespace differences.\n  const current = (cell.value || \"\").trim();\n  if (current !== oldText) {\n    console.log(\n      `Warning: cell (${row}, ${col}) expected \"${oldText}\" but found \"${current}\".`\n    );\n  }\n  cell.value = newText;\n}\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the single table, by\n# position (row, col), matching the diff's old -> new text per cell.\n# Using position rather than literal text Find/Replace avoids ambiguity,\n# since one new value (\"27\u00f75=5, 2\") duplicates an old value used\n# elsewhere in the table.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; Old = \"33\u00f73=11, 0\"; New = \"93\u00f76=15, 3\" },\n    @{ Row = 1; Col = 2; Old = \"32\u00f74=8, 0\"; New = \"85\u00f73=28, 1\" },\n    @{ Row = 1; Col = 3; Old = \"17\u00f75=3, 2\"; New = \"40\u00f77=5, 5\" },\n    @{ Row = 1; Col = 4; Old = \"31\u00f79=3, 4\"; New = \"71\u00f79=7, 8\" },\n    @{ Row = 1; Col = 5; Old = \"90\u00f72=45, 0\"; New = \"28\u00f75=5, 3\" },\n    @{ Row = 5; Col = 1; Old = \"85\u00f74=21, 1\"; New = \"19\u00f72=9, 1\" },\n    @{ Row = 5; Col = 2; Old = \"31\u00f75=6, 1\"; New = \"90\u00f75=18, 0\" },\n    @{ Row = 5; Col = 3; Old = \"61\u00f72=30, 1\"; New = \"88\u00f78=11, 0\" },\n    @{ Row = 5; Col = 4; Old = \"26\u00f79=2, 8\"; New = \"82\u00f77=11, 5\" },\n    @{ Row = 5; Col = 5; Old = \"98\u00f79=10, 8\"; New = \"51\u00f73=17, 0\" },\n    @{ Row = 9; Col = 1; Old = \"99\u00f74=24, 3\"; New = \"81\u00f74=20, 1\" },\n    @{ Row = 9; Col = 2; Old = \"27\u00f75=5, 2\"; New = \"45\u00f78=5, 5\" },\n    @{ Row = 9; Col = 3; Old = \"97\u00f78=12, 1\"; New = \"16\u00f72=8, 0\" },\n    @{ Row = 9; Col = 4; Old = \"16\u00f77=2, 2\"; New = \"58\u00f72=29, 0\" },\n    @{ Row = 9; Col = 5; Old = \"46\u00f72=23, 0\"; New = \"22\u00f78=2, 6\" },\n    @{ Row = 13; Col = 1; Old = \"22\u00f73=7, 1\"; New = \"75\u00f74=18, 3\" },\n    @{ Row = 13; Col = 2; Old = \"84\u00f72=42, 0\"; New = \"57\u00f74=14, 1\" },\n    @{ Row = 13; Col = 3; Old = \"62\u00f73=20, 2\"; New = \"74\u00f72=37, 0\" },\n    @{ Row = 13; Col = 4; Old = \"35\u00f73=11, 2\"; New = \"27\u00f75=5, 2\" },\n    @{ Row = 13; Col = 5; Old = \"43\u00f73=14, 1\"; New = \"18\u00f75=3, 3\" },\n    @{ Row = 17; Col = 1; Old = \"70\u00f72=35, 0\"; New = \"48\u00f72=24, 0\" },\n    @{ Row = 17; Col = 2; Old = \"49\u00f73=16, 1\"; New = \"15\u00f76=2, 3\" },\n    @{ Row = 17; Col = 3; Old = \"56\u00f75=11, 1\"; New = \"88\u00f78=11, 0\" },\n    @{ Row = 17; Col = 4; Old = \"85\u00f77=12, 1\"; New = \"45\u00f72=22, 1\" },\n    @{ Row = 17; Col = 5; Old = \"30\u00f74=7, 2\"; New = \"83\u00f75=16, 3\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n\n    # Sanity-check against the expected original text; Cell.Range.Text\n    # includes the trailing cell-mark characters (CR + BEL), so trim\n    # those off before comparing. The write itself always applies the\n    # intended new value regardless, so the script is resilient to\n    # incidental differences.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $r.Old) {\n        Write-Output \"Warning: cell ($($r.Row), $($r.Col)) expected '$($r.Old)' but found '$current'.\"\n    }\n\n    $cell.Range.Text = $r.New\n}\n"}
